# Trading update: 2026-02-18 10:50:10
#
# - Summary: Total Trades 0 -> 1
# - Strategy Status: MarketMaking row Trades 0 -> 1
# - All Trades: trade #1 (MarketMaking) closes out; two new open trades are
#   logged (#2 momentum, #3 MarketMaking)
# - The old "MarketMaking" per-strategy log sheet is renamed "momentum" and
#   now shows only the momentum trade (#2)
# - A brand new "MarketMaking" per-strategy log sheet is created after it,
#   showing only the MarketMaking trades (#1 and #3)

$wb = $excel.ActiveWorkbook

function Set-TradeRow($ws, $row, $tradeNum, $date, $time, $strategy, $side,
    $entryPrice, $exitPrice, $status, $pnlPct, $pnlDollar, $capitalAfter,
    $entrySlippage, $exitSlippage, $confidence, $entryReason, $exitReason, $duration) {

    $ws.Cells.Item($row, 1).Value = $tradeNum

    # Date/time-looking text must be forced to text so Excel doesn't coerce
    # it into a date/time serial number -- a leading apostrophe is the
    # standard Excel way of doing that.
    $ws.Cells.Item($row, 2).Value = "'" + $date
    $ws.Cells.Item($row, 3).Value = "'" + $time

    $ws.Cells.Item($row, 4).Value = $strategy
    $ws.Cells.Item($row, 5).Value = $side
    $ws.Cells.Item($row, 6).Value = $entryPrice

    if ($null -eq $exitPrice) {
        $ws.Cells.Item($row, 7).Value = ""
    } else {
        $ws.Cells.Item($row, 7).Value = $exitPrice
    }

    $ws.Cells.Item($row, 8).Value = $status
    $ws.Cells.Item($row, 9).Value = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlippage
    $ws.Cells.Item($row, 13).Value = $exitSlippage
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason

    if ([string]::IsNullOrEmpty($exitReason)) {
        $ws.Cells.Item($row, 16).Value = ""
    } else {
        $ws.Cells.Item($row, 16).Value = $exitReason
    }

    $ws.Cells.Item($row, 17).Value = $duration
}

function Set-HeaderRow($ws) {
    $headers = @("Trade #", "Date", "Time", "Strategy", "Side", "Entry Price",
                 "Exit Price", "Status", "P&L %", "P&L $", "Capital After",
                 "Entry Slippage (bps)", "Exit Slippage (bps)", "Confidence",
                 "Entry Reason", "Exit Reason", "Duration (min)")
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    }
}

# ---------------------------------------------------------------------
# Summary: Total Trades 0 -> 1
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 1

# ---------------------------------------------------------------------
# Strategy Status: MarketMaking row -> Trades 0 -> 1
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 1

# ---------------------------------------------------------------------
# All Trades: close out trade #1, append trades #2 and #3
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

Set-TradeRow $allTrades 2 1 "2026-02-18" "10:49:41" "MarketMaking" "UP" `
    0.01 0.01 "CLOSED" 0 0 100 0 0 0.6 "Normal spread capture: 138 bps" "early_exit" 0.09

Set-TradeRow $allTrades 3 2 "2026-02-18" "10:50:02" "momentum" "DOWN" `
    0.47 $null "OPEN" 0 0 100 0 0 0.9 "Downward momentum: -31.724% over 5 samples" $null 0

Set-TradeRow $allTrades 4 3 "2026-02-18" "10:50:02" "MarketMaking" "DOWN" `
    0.479592 $null "OPEN" 0 0 100 0 0 0.65 "Wide spread capture: 202 bps vs avg 151 bps" $null 0

# ---------------------------------------------------------------------
# The old per-strategy "MarketMaking" sheet becomes the per-strategy
# "momentum" log -- rename it in place (keeps its sheetId/rel id) and
# rewrite its single data row to the momentum trade (#2).
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("MarketMaking")
$momentum.Name = "momentum"

Set-TradeRow $momentum 2 2 "2026-02-18" "10:50:02" "momentum" "DOWN" `
    0.47 $null "OPEN" 0 0 100 0 0 0.9 "Downward momentum: -31.724% over 5 samples" $null 0

# ---------------------------------------------------------------------
# Create a brand-new per-strategy "MarketMaking" log sheet right after
# "momentum", with its own two trades (#1 and #3).
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $momentum)
$marketMaking.Name = "MarketMaking"

Set-HeaderRow $marketMaking

Set-TradeRow $marketMaking 2 1 "2026-02-18" "10:49:41" "MarketMaking" "UP" `
    0.01 0.01 "CLOSED" 0 0 100 0 0 0.6 "Normal spread capture: 138 bps" "early_exit" 0.09

Set-TradeRow $marketMaking 3 3 "2026-02-18" "10:50:02" "MarketMaking" "DOWN" `
    0.479592 $null "OPEN" 0 0 100 0 0 0.65 "Wide spread capture: 202 bps vs avg 151 bps" $null 0

$summary.Select()
